$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, border, centered/top aligned) used by the
# rest of row 1 by copying H1's format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-64
$data = @(
 @(5,6),@(8,8),@(6,6),@(5,6),@(5,5),@(7,7),@(7,8),@(7,7),@(6,7),@(9,9),
 @(10,10),@(5,6),@(7,8),@(1,1),@(1,2),@(8,8),@(5,5),@(5,6),@(7,8),@(9,9),
 @(6,6),@(6,6),@(5,6),@(8,8),@(8,8),@(8,8),@(6,6),@(5,6),@(6,6),@(8,8),
 @(5,5),@(8,8),@(8,8),@(9,9),@(9,9),@(6,7),@(7,7),@(6,7),@(8,8),@(6,6),
 @(7,8),@(8,8),@(11,11),@(7,7),@(8,8),@(8,8),@(8,8),@(7,7),@(8,8),@(8,8),
 @(8,8),@(7,7),@(6,6),@(5,6),@(8,8),@(10,10),@(8,8),@(8,8),@(5,5),@(7,7),
 @(6,6),@(5,5),@(7,7)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
